# Add new column 'Servised by' to Card18 (by admin)
#
# The sheet currently spans A1:N12. We append a new column O:
#   - O1 gets the header "Servised by", formatted like the other
#     header cells (bold/centered/bordered - same style as N1).
#   - O2:O12 are created as (empty) text cells, mirroring the other
#     "new" empty cells already present in the sheet.
#   - N2:N12, which previously held an empty inline string, become
#     explicit "nan" text values (matching the rest of the column's
#     existing "nan" placeholders).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")
$ws.Activate()

# Header cell for the new column.
$ws.Range("O1").Value = "Servised by"

# Copy the header formatting (bold font, border, centered alignment)
# from the existing "Correction" header so the new header matches.
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)  # xlPasteFormats

for ($r = 2; $r -le 12; $r++) {
    # N column previously had an empty placeholder cell; fill it with
    # the same "nan" placeholder used throughout the rest of the row.
    $ws.Cells.Item($r, 14).Value = "nan"

    # Materialize an empty text cell in the new O column (leading
    # apostrophe forces a blank text entry instead of leaving the
    # cell completely unset), then strip the quote-prefix formatting
    # it implies so the cell keeps the default style.
    $ws.Cells.Item($r, 15).Value = "'"
    $ws.Cells.Item($r, 15).Style = "Normal"
}
